$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D (Price) contains numeric-looking text (e.g. "1.001", "26.461.90")
# that Excel would otherwise auto-convert to real numbers. Force the whole
# price column to Text format before writing the new values, then restore
# the default "Normal" style so no stray formatting is left behind.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

$ws.Range("D2").Value = "26.461.90"
$ws.Range("D3").Value = "1.835.82"
$ws.Range("D6").Value = "1.001"
$ws.Range("D8").Value = "0.3000"
$ws.Range("D9").Value = "0.06893"
$ws.Range("D10").Value = "17.61"
$ws.Range("D11").Value = "1.849.43"
$ws.Range("D12").Value = "0.7354"
$ws.Range("D13").Value = "0.07164"
$ws.Range("D14").Value = "88.93"
$ws.Range("D19").Value = "0.000007909"
$ws.Range("D20").Value = "26.490.51"
$ws.Range("D21").Value = "2.079.34"
$ws.Range("D22").Value = "4.597"
$ws.Range("D24").Value = "9.204"
$ws.Range("D25").Value = "142.93"
$ws.Range("D26").Value = "2.182"
$ws.Range("D27").Value = "1.720"
$ws.Range("D28").Value = "17.00"
$ws.Range("D29").Value = "110.91"
$ws.Range("D30").Value = "4.249"
$ws.Range("D31").Value = "0.08840"
$ws.Range("D32").Value = "4.041"
$ws.Range("D34").Value = "2.919"
$ws.Range("D35").Value = "0.7295"
$ws.Range("D39").Value = "0.01714"
$ws.Range("D40").Value = "0.4711"
$ws.Range("D41").Value = "0.9022"
$ws.Range("D42").Value = "108.07"
$ws.Range("D43").Value = "5.895"
$ws.Range("D45").Value = "7.394"
$ws.Range("D46").Value = "0.1250"
$ws.Range("D47").Value = "9.010"
$ws.Range("D48").Value = "0.4070"
$ws.Range("D49").Value = "34.81"
$ws.Range("D50").Value = "0.8930"
$ws.Range("D51").Value = "0.05768"

$priceRange.Style = "Normal"

# Column E (Volume(1h)) values are already padded percentage strings and
# stay text on assignment, so they can be set directly.
$ws.Range("E2").Value = "  -0.55%  "
$ws.Range("E4").Value = "  -0.16%  "
$ws.Range("E5").Value = "  -1.00%  "
$ws.Range("E6").Value = "  -0.11%  "
$ws.Range("E7").Value = "  +2.15%  "
$ws.Range("E8").Value = "  -7.57%  "
$ws.Range("E9").Value = "  +1.41%  "
$ws.Range("E10").Value = "  -6.95%  "
$ws.Range("E11").Value = "  -0.24%  "
$ws.Range("E12").Value = "  -6.05%  "
$ws.Range("E13").Value = "  -7.94%  "
$ws.Range("E14").Value = "  +0.37%  "
$ws.Range("E15").Value = "  -0.61%  "
$ws.Range("E16").Value = "  -0.11%  "
$ws.Range("E17").Value = "  -1.04%  "
$ws.Range("E18").Value = "  -0.11%  "
$ws.Range("E19").Value = "  -0.79%  "
$ws.Range("E20").Value = "  -0.58%  "
$ws.Range("E21").Value = "  -0.73%  "
$ws.Range("E22").Value = "  -0.92%  "
$ws.Range("E23").Value = "  -0.51%  "
$ws.Range("E24").Value = "  -2.82%  "
$ws.Range("E25").Value = "  -0.41%  "
$ws.Range("E26").Value = "  +0.60%  "
$ws.Range("E27").Value = "  +1.65%  "
$ws.Range("E28").Value = "  -0.26%  "
$ws.Range("E29").Value = "  -0.94%  "
$ws.Range("E30").Value = "  +1.49%  "
$ws.Range("E31").Value = "  +1.39%  "
$ws.Range("E32").Value = "  -1.79%  "
$ws.Range("E33").Value = "  -0.55%  "
$ws.Range("E34").Value = "  +1.70%  "
$ws.Range("E35").Value = "  +1.30%  "
$ws.Range("E36").Value = "  +0.06%  "
$ws.Range("E37").Value = "  -0.78%  "
$ws.Range("E38").Value = "  -0.49%  "
$ws.Range("E39").Value = "  -4.52%  "
$ws.Range("E40").Value = "  -3.33%  "
$ws.Range("E41").Value = "  +0.00%  "
$ws.Range("E42").Value = "  -2.77%  "
$ws.Range("E43").Value = "  -1.14%  "
$ws.Range("E44").Value = "  -0.12%  "
$ws.Range("E45").Value = "  -3.61%  "
$ws.Range("E46").Value = "  +1.30%  "
$ws.Range("E47").Value = "  -0.50%  "
$ws.Range("E48").Value = "  -2.83%  "
$ws.Range("E49").Value = "  -0.86%  "
$ws.Range("E50").Value = "  +0.41%  "
$ws.Range("E51").Value = "  -2.12%  "

